$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(1)
$para = $sh.TextFrame.TextRange.Paragraphs(1)
$para.Text = "X"
$para.Text = "Below section-level"
